$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update team-specific time transition matrix values
# Row 2
$ws.Cells.Item(2, 2).Value = 0.1582089552238806
$ws.Cells.Item(2, 3).Value = 0.6119402985074627
$ws.Cells.Item(2, 10).Value = 0.0208955223880597
$ws.Cells.Item(2, 16).Value = 0.1223880597014925
$ws.Cells.Item(2, 19).Value = 0.08656716417910448

# Row 3
$ws.Cells.Item(3, 2).Value = 0.004587155963302753
$ws.Cells.Item(3, 3).Value = 0.04128440366972477
$ws.Cells.Item(3, 10).Value = 0.04587155963302753
$ws.Cells.Item(3, 16).Value = 0.7201834862385321
$ws.Cells.Item(3, 19).Value = 0.1880733944954129

# Row 4
$ws.Cells.Item(4, 10).Value = 0.1
$ws.Cells.Item(4, 16).Value = 0.725
$ws.Cells.Item(4, 19).Value = 0.175

# Row 5
$ws.Cells.Item(5, 16).Value = 0.6
$ws.Cells.Item(5, 19).Value = 0.4

# Row 6
$ws.Cells.Item(6, 2).Value = 0.06589147286821706
$ws.Cells.Item(6, 6).Value = 0.07364341085271318
$ws.Cells.Item(6, 10).Value = 0.2596899224806202
$ws.Cells.Item(6, 15).Value = 0.01937984496124031
$ws.Cells.Item(6, 17).Value = 0.1550387596899225
$ws.Cells.Item(6, 18).Value = 0.05426356589147287
$ws.Cells.Item(6, 19).Value = 0.3720930232558139

# Row 7
$ws.Cells.Item(7, 2).Value = 0.09359605911330049
$ws.Cells.Item(7, 4).Value = 0.03448275862068965
$ws.Cells.Item(7, 6).Value = 0.04926108374384237
$ws.Cells.Item(7, 10).Value = 0.1970443349753695
$ws.Cells.Item(7, 15).Value = 0.01970443349753695
$ws.Cells.Item(7, 17).Value = 0.1231527093596059
$ws.Cells.Item(7, 18).Value = 0.1330049261083744
$ws.Cells.Item(7, 19).Value = 0.3497536945812808

# Row 8
$ws.Cells.Item(8, 2).Value = 0.08812260536398467
$ws.Cells.Item(8, 4).Value = 0.01915708812260536
$ws.Cells.Item(8, 5).Value = 0.003831417624521073
$ws.Cells.Item(8, 6).Value = 0.04597701149425287
$ws.Cells.Item(8, 10).Value = 0.1283524904214559
$ws.Cells.Item(8, 15).Value = 0.0210727969348659
$ws.Cells.Item(8, 17).Value = 0.181992337164751
$ws.Cells.Item(8, 18).Value = 0.1053639846743295
$ws.Cells.Item(8, 19).Value = 0.4061302681992337

# Row 9
$ws.Cells.Item(9, 2).Value = 0.07860262008733625
$ws.Cells.Item(9, 4).Value = 0.01310043668122271
$ws.Cells.Item(9, 6).Value = 0.05676855895196507
$ws.Cells.Item(9, 10).Value = 0.148471615720524
$ws.Cells.Item(9, 15).Value = 0.008733624454148471
$ws.Cells.Item(9, 17).Value = 0.1790393013100437
$ws.Cells.Item(9, 18).Value = 0.1397379912663755
$ws.Cells.Item(9, 19).Value = 0.3755458515283843

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1140808344198175
$ws.Cells.Item(10, 4).Value = 0.01434159061277705
$ws.Cells.Item(10, 5).Value = 0.002607561929595828
$ws.Cells.Item(10, 6).Value = 0.06910039113428944
$ws.Cells.Item(10, 10).Value = 0.1544980443285528
$ws.Cells.Item(10, 15).Value = 0.01499348109517601
$ws.Cells.Item(10, 17).Value = 0.1799217731421121
$ws.Cells.Item(10, 18).Value = 0.08865710560625815
$ws.Cells.Item(10, 19).Value = 0.3617992177314211

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1474358974358974
$ws.Cells.Item(11, 10).Value = 0.08653846153846154
$ws.Cells.Item(11, 11).Value = 0.1955128205128205
$ws.Cells.Item(11, 12).Value = 0.5480769230769231
$ws.Cells.Item(11, 19).Value = 0.02243589743589744

# Row 12
$ws.Cells.Item(12, 6).Value = 0.005434782608695652
$ws.Cells.Item(12, 7).Value = 0.7065217391304348
$ws.Cells.Item(12, 10).Value = 0.1793478260869565
$ws.Cells.Item(12, 11).Value = 0.02717391304347826
$ws.Cells.Item(12, 12).Value = 0.05978260869565218
$ws.Cells.Item(12, 19).Value = 0.02173913043478261

# Row 13
$ws.Cells.Item(13, 7).Value = 0.575
$ws.Cells.Item(13, 10).Value = 0.375
$ws.Cells.Item(13, 19).Value = 0.05

# Row 14
$ws.Cells.Item(14, 7).Value = 0.9
$ws.Cells.Item(14, 10).Value = 0.1

# Row 15
$ws.Cells.Item(15, 6).Value = 0.0282258064516129
$ws.Cells.Item(15, 8).Value = 0.1330645161290323
$ws.Cells.Item(15, 9).Value = 0.07258064516129033
$ws.Cells.Item(15, 10).Value = 0.4112903225806452
$ws.Cells.Item(15, 11).Value = 0.0564516129032258
$ws.Cells.Item(15, 15).Value = 0.04032258064516129
$ws.Cells.Item(15, 19).Value = 0.2580645161290323

# Row 16
$ws.Cells.Item(16, 6).Value = 0.03167420814479638
$ws.Cells.Item(16, 8).Value = 0.1583710407239819
$ws.Cells.Item(16, 9).Value = 0.1040723981900453
$ws.Cells.Item(16, 10).Value = 0.3981900452488688
$ws.Cells.Item(16, 11).Value = 0.08597285067873303
$ws.Cells.Item(16, 13).Value = 0.01357466063348416
$ws.Cells.Item(16, 14).Value = 0.009049773755656109
$ws.Cells.Item(16, 15).Value = 0.04524886877828054
$ws.Cells.Item(16, 19).Value = 0.1538461538461539

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02564102564102564
$ws.Cells.Item(17, 8).Value = 0.1837606837606838
$ws.Cells.Item(17, 9).Value = 0.08974358974358974
$ws.Cells.Item(17, 10).Value = 0.4230769230769231
$ws.Cells.Item(17, 11).Value = 0.07478632478632478
$ws.Cells.Item(17, 13).Value = 0.01495726495726496
$ws.Cells.Item(17, 15).Value = 0.05982905982905983
$ws.Cells.Item(17, 19).Value = 0.1282051282051282

# Row 18
$ws.Cells.Item(18, 6).Value = 0.01520912547528517
$ws.Cells.Item(18, 8).Value = 0.1711026615969582
$ws.Cells.Item(18, 9).Value = 0.1140684410646388
$ws.Cells.Item(18, 10).Value = 0.403041825095057
$ws.Cells.Item(18, 11).Value = 0.1026615969581749
$ws.Cells.Item(18, 13).Value = 0.007604562737642586
$ws.Cells.Item(18, 15).Value = 0.05703422053231939
$ws.Cells.Item(18, 19).Value = 0.1292775665399239

# Row 19
$ws.Cells.Item(19, 6).Value = 0.01929703652653342
$ws.Cells.Item(19, 8).Value = 0.217091660923501
$ws.Cells.Item(19, 9).Value = 0.0771881461061337
$ws.Cells.Item(19, 10).Value = 0.3659545141281875
$ws.Cells.Item(19, 11).Value = 0.1026878015161957
$ws.Cells.Item(19, 13).Value = 0.01998621640248105
$ws.Cells.Item(19, 14).Value = 0.006202618883528601
$ws.Cells.Item(19, 15).Value = 0.07649896623018608
$ws.Cells.Item(19, 19).Value = 0.1150930392832529
